$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2.092578887939453, -8.521716117858887, 1.690641283988953),
    @(2.046772480010986, -4.874783515930176, 0.8943560719490051),
    @(4.356798648834229, -5.594369411468506, -4.091081619262695),
    @(-5.464406967163086, 3.147723197937012, 2.375073671340942),
    @(-6.766693115234375, 7.695337295532227, 7.492071628570557),
    @(3.930692672729492, 9.625063896179199, -5.505752563476562),
    @(3.03800106048584, -1.461676001548767, 0.6908905506134033),
    @(11.57449817657471, -5.450558662414551, -0.0244345031678676),
    @(-4.411925792694092, 2.68752908706665, 3.053647041320801),
    @(-4.802345275878906, 4.841493606567383, 2.432597875595093),
    @(-2.648380517959595, 16.29841232299805, -2.319014072418213),
    @(-0.7825698852539062, -1.549027681350708, 1.645367503166199),
    @(9.456219673156738, -6.493985176086426, -6.240252494812012),
    @(3.910985231399536, 0.8456867933273315, -5.051949977874756),
    @(-2.794321775436401, -2.126400947570801, -1.946171641349792),
    @(-14.5798749923706, -11.00378227233887, 6.221210956573486),
    @(5.833254814147949, -11.42988777160644, 6.618554592132568),
    @(0.6635265946388245, -4.405001640319824, 1.273590207099915),
    @(5.506218910217285, 16.36765480041504, -1.815144062042236),
    @(-1.695501565933228, -2.446512937545776, 2.133258581161499),
    @(-4.14081621170044, 5.169595241546631, 1.407280921936035),
    @(-2.688860654830933, 10.50284194946289, -1.326720356941223),
    @(-1.443033814430237, 10.46609020233154, -6.833072185516357),
    @(1.430516958236694, -2.983938932418823, -5.8599534034729),
    @(3.719237804412842, 1.893374443054199, -7.143064022064209),
    @(-0.8805742263793945, -0.7974836230278015, 1.437640905380249),
    @(-4.217514991760254, 1.791641712188721, 6.048105716705322),
    @(-1.25075364112854, 12.15346908569336, -1.156278014183044),
    @(-1.032907009124756, 9.82000732421875, -6.767558574676514),
    @(-3.240667581558228, -12.29328536987305, -4.838364601135254)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
}